$d = $word.ActiveDocument

function Replace-ParagraphText($oldText, $newText) {
    $r = $d.Content
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r.Text = $newText
}

# 1) Replace the first "Detailed Feedback" bullet's text.
Replace-ParagraphText `
    "The interviewee, Goh Yi Xian, demonstrates some relevant technical skills in deep learning and machine learning, as evidenced by their experience in medical image processing and research activities in IVIF. However, there are notable gaps in meeting the job requirements for an LLM Engineer, particularly in the areas of NLP, large language model development, and deployment." `
    "The interviewee, Goh Yi Xian, shows potential with relevant work experience in machine learning and research activities in image processing. However, there are significant gaps in meeting the technical skills and experience required for the LLM Engineer role."

# 2) Replace the second bullet's text.
Replace-ParagraphText `
    "While Goh Yi Xian shows potential in problem-solving and innovation, their response in the interview lacked clarity, relevance, and depth in addressing key aspects of LLM development, data management, and model evaluation." `
    "While Goh Yi Xian has experience in deep learning and image processing, there is a lack of demonstrated expertise in developing and deploying large language models (LLMs) for NLP tasks. The interviewee's response in the interview did not address key aspects of model development, data management, infrastructure deployment, and model evaluation as outlined in the job requirements."

# 3) Replace the third bullet's text.
Replace-ParagraphText `
    "The interviewee's confidence score of 100 did not translate effectively into the content of their responses, indicating a disconnect between presentation and substance." `
    "The educational background and soft skills of Goh Yi Xian align well with the job description, showcasing strong problem-solving abilities and a proactive mindset. However, the lack of specific experience in NLP and transformer-based models like BERT or GPT is a notable gap."

# 4) Remove the fourth bullet paragraph entirely (personality traits bullet).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Goh Yi Xian's personality traits of being detail-oriented*") {
        $p.Range.Delete()
        break
    }
}

# 5) Replace the fifth (now final) bullet's text.
Replace-ParagraphText `
    "Overall, Goh Yi Xian's performance in this interview falls below expectations for the LLM Engineer role. There is a need for further development in NLP, deep learning, and model deployment skills to be a strong fit for the position." `
    "In summary, while Goh Yi Xian shows promise in the AI field, further development in NLP, deep learning frameworks, and LLM deployment is needed to be a suitable match for the LLM Engineer role."
